$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 9562.15
$ws.Range("B16").Value = 10028.469999999999
$ws.Range("C16").Value = 17.2
$ws.Range("D16").Value = 18
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = 4.6500000000000004
$ws.Range("G15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 42626.545613425929
$ws.Range("H16").Value = $false
